# Updates SettingApp.xlsx:
#  - Splits the "RutasArcchivos" table out of the "Setting" sheet into its
#    own new "Rutas" sheet (placed between "Setting" and "Tornilleria").
#  - Extends the "Inventor" table on "Setting" with four new rows/columns
#    names: "CTDADes", "Material", "CTDAD de elementos" and "Proveedor",
#    while removing the old standalone "CTDAD" entry.

$wb = $excel.ActiveWorkbook
$wsSetting = $wb.Worksheets.Item("Setting")
$wsTornilleria = $wb.Worksheets.Item("Tornilleria")

# ---------------------------------------------------------------------------
# 1) Pull the RutasArcchivos table off of "Setting" (columns C:D) so it can
#    be rebuilt on its own sheet.
# ---------------------------------------------------------------------------
$rutasLo = $wsSetting.ListObjects.Item("RutasArcchivos")
$rutasLo.Delete()
$wsSetting.Range("C1:D9").Clear()

# ---------------------------------------------------------------------------
# 2) Rewrite the Inventor column list (column A) on "Setting" with the new
#    12-row layout (adds Material / CTDAD de elementos / Proveedor / CTDADes,
#    drops the standalone "CTDAD" row). The brand-new names are written in
#    this specific order so the workbook's shared-string table grows the same
#    way it did in the authored edit (new uniques appended in entry order,
#    not sheet/row order).
# ---------------------------------------------------------------------------
$wsSetting.Cells.Item(1, 1).Value = "ColumnasExcelInventor"
$wsSetting.Cells.Item(2, 1).Value = "Elemento"
$wsSetting.Cells.Item(4, 1).Value = "Nº de pieza"
$wsSetting.Cells.Item(5, 1).Value = "Descripción"
$wsSetting.Cells.Item(6, 1).Value = "CTDAD de unidades"
$wsSetting.Cells.Item(9, 1).Value = "Masa"
$wsSetting.Cells.Item(10, 1).Value = "Nombre de archivo"
$wsSetting.Cells.Item(12, 1).Value = "Tipo de componente"

$wsSetting.Cells.Item(7, 1).Value = "Material"
$wsSetting.Cells.Item(8, 1).Value = "CTDAD de elementos"
$wsSetting.Cells.Item(11, 1).Value = "Proveedor"
$wsSetting.Cells.Item(3, 1).Value = "CTDADes"

$inventorLo = $wsSetting.ListObjects.Item("Inventor")
$inventorLo.Resize($wsSetting.Range("A1:A12"))

$wsSetting.Range("A1:D1048576").EntireColumn.AutoFit()
$wsSetting.Application.ActiveWindow.Zoom = 100
$wsSetting.Range("C5").Select()

# ---------------------------------------------------------------------------
# 3) Insert the new "Rutas" sheet right after "Setting" and build the
#    RutasArcchivos table there (columns A:B).
# ---------------------------------------------------------------------------
$wsRutas = $wb.Worksheets.Add($null, $wsSetting)
$wsRutas.Name = "Rutas"

$wsRutas.Range("A1").Value = "RutasArchivo"
$wsRutas.Range("B1").Value = "Direccion"
$wsRutas.Range("A2").Value = "Inventor"
$wsRutas.Range("B2").Value = "C:\Users\myb19\Desktop\ID011 LM Horno Estructurado.xlsx"
$wsRutas.Range("A3").Value = "Maestro"
$wsRutas.Range("B3").Value = "C:\Users\myb19\Desktop\Maestro Productos Molecor TECH.xlsx"
$wsRutas.Range("A4").Value = "Navision"
$wsRutas.Range("A5").Value = "Sap"

$rutasNewLo = $wsRutas.ListObjects.Add(1, $wsRutas.Range("A1:B5"), $null, 1)
$rutasNewLo.Name = "RutasArcchivos"

$wsRutas.Columns.Item(1).ColumnWidth = 14.85546875
$wsRutas.Columns.Item(2).ColumnWidth = 62
$wsRutas.Range("B12").Select()

# ---------------------------------------------------------------------------
# 4) Leave "Setting" as the active/selected sheet, matching the saved state.
# ---------------------------------------------------------------------------
$wsSetting.Activate()
$wsSetting.Range("C5").Select()
